$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.3642143333333334
$ws.Cells.Item(2, 8).Value = 1.092643
$ws.Cells.Item(2, 9).Value = 0.4800482050304226
$ws.Cells.Item(2, 10).Value = 0.4800482050304224
$ws.Cells.Item(2, 13).Value = 19.77408333333333
$ws.Cells.Item(2, 14).Value = 59.32225
$ws.Cells.Item(2, 15).Value = 0.3380388258879848
$ws.Cells.Item(2, 16).Value = 0.339186328349942
$ws.Cells.Item(2, 17).Value = 7.202004578527778
$ws.Cells.Item(2, 18).Value = 64.81804120675001
$ws.Cells.Item(2, 19).Value = 0.1622749315981186
$ws.Cells.Item(2, 20).Value = 0.1628257880952491

# Row 3
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.3642143333333334
$ws.Cells.Item(3, 8).Value = 1.092643
$ws.Cells.Item(3, 9).Value = 0.4800482050304226
$ws.Cells.Item(3, 10).Value = 0.4800482050304224
$ws.Cells.Item(3, 15).Value = 0.3069959581674471
$ws.Cells.Item(3, 16).Value = 0.3080380828904952
$ws.Cells.Item(3, 17).Value = 6.54062825624689
$ws.Cells.Item(3, 18).Value = 58.86565430622201
$ws.Cells.Item(3, 19).Value = 0.1473728586698776
$ws.Cells.Item(3, 20).Value = 0.1478731287725947

# Row 4
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.3642143333333334
$ws.Cells.Item(4, 8).Value = 1.092643
$ws.Cells.Item(4, 9).Value = 0.4800482050304226
$ws.Cells.Item(4, 10).Value = 0.4800482050304224
$ws.Cells.Item(4, 13).Value = 12.46730333333333
$ws.Cells.Item(4, 14).Value = 37.40191
$ws.Cells.Item(4, 15).Value = 0.2131290998296268
$ws.Cells.Item(4, 16).Value = 0.2138525852639604
$ws.Cells.Item(4, 17).Value = 4.540770572014445
$ws.Cells.Item(4, 18).Value = 40.86693514813
$ws.Cells.Item(4, 19).Value = 0.1023122418129621
$ws.Cells.Item(4, 20).Value = 0.1026595496970796

# Row 5
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.3642143333333334
$ws.Cells.Item(5, 8).Value = 1.092643
$ws.Cells.Item(5, 9).Value = 0.4800482050304226
$ws.Cells.Item(5, 10).Value = 0.4800482050304224
$ws.Cells.Item(5, 13).Value = 0.593699
$ws.Cells.Item(5, 14).Value = 1.187398
$ws.Cells.Item(5, 15).Value = 0.01014931056513554
$ws.Cells.Item(5, 16).Value = 0.006789175527058808
$ws.Cells.Item(5, 17).Value = 0.2162336854856667
$ws.Cells.Item(5, 18).Value = 1.297402112914
$ws.Cells.Item(5, 19).Value = 0.004872158319089622
$ws.Cells.Item(5, 20).Value = 0.003259131525401053

# Row 6
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.3642143333333334
$ws.Cells.Item(6, 8).Value = 1.092643
$ws.Cells.Item(6, 9).Value = 0.4800482050304226
$ws.Cells.Item(6, 10).Value = 0.4800482050304224
$ws.Cells.Item(6, 13).Value = 7.703215333333333
$ws.Cells.Item(6, 14).Value = 23.109646
$ws.Cells.Item(6, 15).Value = 0.1316868055498057
$ws.Cells.Item(6, 16).Value = 0.1321338279685434
$ws.Cells.Item(6, 17).Value = 2.805621437153111
$ws.Cells.Item(6, 18).Value = 25.250592934378
$ws.Cells.Item(6, 19).Value = 0.06321601463037453
$ws.Cells.Item(6, 20).Value = 0.06343060694009788

# Row 7
$ws.Cells.Item(7, 7).Value = 0.3944893333333333
$ws.Cells.Item(7, 8).Value = 1.183468
$ws.Cells.Item(7, 9).Value = 0.5199517949695774
$ws.Cells.Item(7, 10).Value = 0.5199517949695774
$ws.Cells.Item(7, 13).Value = 19.77408333333333
$ws.Cells.Item(7, 14).Value = 59.32225
$ws.Cells.Item(7, 15).Value = 0.3380388258879848
$ws.Cells.Item(7, 16).Value = 0.339186328349942
$ws.Cells.Item(7, 17).Value = 7.800664951444443
$ws.Cells.Item(7, 18).Value = 70.205984563
$ws.Cells.Item(7, 19).Value = 0.1757638942898662
$ws.Cells.Item(7, 20).Value = 0.1763605402546928

# Row 8
$ws.Cells.Item(8, 7).Value = 0.3944893333333333
$ws.Cells.Item(8, 8).Value = 1.183468
$ws.Cells.Item(8, 9).Value = 0.5199517949695774
$ws.Cells.Item(8, 10).Value = 0.5199517949695774
$ws.Cells.Item(8, 15).Value = 0.3069959581674471
$ws.Cells.Item(8, 16).Value = 0.3080380828904952
$ws.Cells.Item(8, 17).Value = 7.084312297030222
$ws.Cells.Item(8, 18).Value = 63.758810673272
$ws.Cells.Item(8, 19).Value = 0.1596230994975694
$ws.Cells.Item(8, 20).Value = 0.1601649541179005

# Row 9
$ws.Cells.Item(9, 7).Value = 0.3944893333333333
$ws.Cells.Item(9, 8).Value = 1.183468
$ws.Cells.Item(9, 9).Value = 0.5199517949695774
$ws.Cells.Item(9, 10).Value = 0.5199517949695774
$ws.Cells.Item(9, 13).Value = 12.46730333333333
$ws.Cells.Item(9, 14).Value = 37.40191
$ws.Cells.Item(9, 15).Value = 0.2131290998296268
$ws.Cells.Item(9, 16).Value = 0.2138525852639604
$ws.Cells.Item(9, 17).Value = 4.91821818043111
$ws.Cells.Item(9, 18).Value = 44.26396362388
$ws.Cells.Item(9, 19).Value = 0.1108168580166647
$ws.Cells.Item(9, 20).Value = 0.1111930355668808

# Row 10
$ws.Cells.Item(10, 7).Value = 0.3944893333333333
$ws.Cells.Item(10, 8).Value = 1.183468
$ws.Cells.Item(10, 9).Value = 0.5199517949695774
$ws.Cells.Item(10, 10).Value = 0.5199517949695774
$ws.Cells.Item(10, 13).Value = 0.593699
$ws.Cells.Item(10, 14).Value = 1.187398
$ws.Cells.Item(10, 15).Value = 0.01014931056513554
$ws.Cells.Item(10, 16).Value = 0.006789175527058808
$ws.Cells.Item(10, 17).Value = 0.2342079227106666
$ws.Cells.Item(10, 18).Value = 1.405247536264
$ws.Cells.Item(10, 19).Value = 0.005277152246045923
$ws.Cells.Item(10, 20).Value = 0.003530044001657754

# Row 11
$ws.Cells.Item(11, 7).Value = 0.3944893333333333
$ws.Cells.Item(11, 8).Value = 1.183468
$ws.Cells.Item(11, 9).Value = 0.5199517949695774
$ws.Cells.Item(11, 10).Value = 0.5199517949695774
$ws.Cells.Item(11, 13).Value = 7.703215333333333
$ws.Cells.Item(11, 14).Value = 23.109646
$ws.Cells.Item(11, 15).Value = 0.1316868055498057
$ws.Cells.Item(11, 16).Value = 0.1321338279685434
$ws.Cells.Item(11, 17).Value = 3.038836281369777
$ws.Cells.Item(11, 18).Value = 27.349526532328
$ws.Cells.Item(11, 19).Value = 0.06847079091943122
$ws.Cells.Item(11, 20).Value = 0.06870322102844548
